$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B2').Value = 'Tue Jun 03 20:37:33 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B3').Value = 'Tue Jun 03 20:38:22 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B4').Value = 'Tue Jun 03 20:39:06 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B5').Value = 'Tue Jun 03 20:39:55 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B6').Value = 'Tue Jun 03 20:40:47 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B7').Value = 'Tue Jun 03 20:41:34 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B8').Value = 'Tue Jun 03 20:42:18 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B9').Value = 'Tue Jun 03 20:43:10 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPC')
$ws.Range('B2').Value = 'Tue Jun 03 19:56:55 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPS')
$ws.Range('B2').Value = 'Tue Jun 03 19:58:23 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPS')
$ws.Range('B3').Value = 'Tue Jun 03 19:59:17 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPS')
$ws.Range('B4').Value = 'Tue Jun 03 20:00:09 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPS')
$ws.Range('B5').Value = 'Tue Jun 03 20:01:01 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPS')
$ws.Range('B6').Value = 'Tue Jun 03 20:01:47 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPS')
$ws.Range('B7').Value = 'Tue Jun 03 20:02:35 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorp')
$ws.Range('B2').Value = 'Tue Jun 03 20:11:50 IST 2025'

$ws = $wb.Worksheets.Item('MRFCorp')
$ws.Range('B2').Value = 'Tue Jun 03 22:41:52 IST 2025'

$ws = $wb.Worksheets.Item('MRFCorp')
$ws.Range('B3').Value = 'Tue Jun 03 22:42:37 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsCorp')
$ws.Range('B2').Value = 'Tue Jun 03 21:08:14 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsCorp')
$ws.Range('B3').Value = 'Tue Jun 03 21:08:54 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsCorp')
$ws.Range('B2').Value = 'Tue Jun 03 21:38:09 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsCorp')
$ws.Range('B3').Value = 'Tue Jun 03 21:38:50 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsPS')
$ws.Range('B2').Value = 'Tue Jun 03 21:10:52 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsPS')
$ws.Range('B3').Value = 'Tue Jun 03 21:11:37 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsPS')
$ws.Range('B2').Value = 'Tue Jun 03 21:41:27 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsPS')
$ws.Range('B3').Value = 'Tue Jun 03 21:42:17 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsCredit')
$ws.Range('B2').Value = 'Tue Jun 03 21:06:59 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsCredit')
$ws.Range('B3').Value = 'Tue Jun 03 21:07:39 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsCredit')
$ws.Range('B2').Value = 'Tue Jun 03 21:36:50 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsCredit')
$ws.Range('A3').Value = 'Fail'
$ws.Range('B3').Value = 'Tue Jun 03 21:37:31 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsPC')
$ws.Range('B2').Value = 'Tue Jun 03 21:09:32 IST 2025'

$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsPC')
$ws.Range('B3').Value = 'Tue Jun 03 21:10:12 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsPC')
$ws.Range('B2').Value = 'Tue Jun 03 21:39:38 IST 2025'

$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsPC')
$ws.Range('B3').Value = 'Tue Jun 03 21:40:32 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpSCF')
$ws.Range('B2').Value = 'Tue Jun 03 20:07:33 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpSCF')
$ws.Range('B3').Value = 'Tue Jun 03 20:08:41 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpSCF')
$ws.Range('B4').Value = 'Tue Jun 03 20:09:47 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpSCF')
$ws.Range('B5').Value = 'Tue Jun 03 20:10:50 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditSCF')
$ws.Range('B2').Value = 'Tue Jun 03 20:16:44 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditSCF')
$ws.Range('B3').Value = 'Tue Jun 03 20:17:52 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditSCF')
$ws.Range('B4').Value = 'Tue Jun 03 20:18:54 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditSCF')
$ws.Range('B5').Value = 'Tue Jun 03 20:19:55 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckSCF')
$ws.Range('B2').Value = 'Tue Jun 03 20:25:09 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckSCF')
$ws.Range('B3').Value = 'Tue Jun 03 20:26:14 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckSCF')
$ws.Range('B4').Value = 'Tue Jun 03 20:27:14 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckSCF')
$ws.Range('B5').Value = 'Tue Jun 03 20:28:20 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalSavingsSCF')
$ws.Range('B2').Value = 'Tue Jun 03 20:29:20 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalSavingsSCF')
$ws.Range('B3').Value = 'Tue Jun 03 20:30:24 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalSavingsSCF')
$ws.Range('B4').Value = 'Tue Jun 03 20:31:27 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalSavingsSCF')
$ws.Range('B5').Value = 'Tue Jun 03 20:32:28 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditDCF')
$ws.Range('B2').Value = 'Tue Jun 03 20:12:35 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditDCF')
$ws.Range('B3').Value = 'Tue Jun 03 20:13:40 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditDCF')
$ws.Range('B4').Value = 'Tue Jun 03 20:14:41 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCreditDCF')
$ws.Range('B5').Value = 'Tue Jun 03 20:15:42 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpDCF')
$ws.Range('B2').Value = 'Tue Jun 03 20:03:20 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpDCF')
$ws.Range('B3').Value = 'Tue Jun 03 20:04:23 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpDCF')
$ws.Range('B4').Value = 'Tue Jun 03 20:05:26 IST 2025'

$ws = $wb.Worksheets.Item('PayNowCorpDCF')
$ws.Range('B5').Value = 'Tue Jun 03 20:06:26 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckDCF')
$ws.Range('B2').Value = 'Tue Jun 03 20:33:27 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckDCF')
$ws.Range('B3').Value = 'Tue Jun 03 20:34:28 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckDCF')
$ws.Range('B4').Value = 'Tue Jun 03 20:35:29 IST 2025'

$ws = $wb.Worksheets.Item('PayNowPersonalCheckDCF')
$ws.Range('B5').Value = 'Tue Jun 03 20:36:31 IST 2025'

$ws = $wb.Worksheets.Item('MaxAmountErrorCC')
$ws.Range('B3').Value = 'Tue Jun 03 22:30:25 IST 2025'

$ws = $wb.Worksheets.Item('MaxAmountErrorCorp')
$ws.Range('B2').Value = 'Tue Jun 03 22:31:09 IST 2025'

$ws = $wb.Worksheets.Item('MaxAmountErrorCorp')
$ws.Range('B3').Value = 'Tue Jun 03 22:31:55 IST 2025'

$ws = $wb.Worksheets.Item('MaxAmountErrorPC')
$ws.Range('B2').Value = 'Tue Jun 03 22:32:40 IST 2025'

$ws = $wb.Worksheets.Item('MaxAmountErrorPC')
$ws.Range('B3').Value = 'Tue Jun 03 22:33:23 IST 2025'

$ws = $wb.Worksheets.Item('MaxAmountErrorPS')
$ws.Range('B2').Value = 'Tue Jun 03 22:34:11 IST 2025'

$ws = $wb.Worksheets.Item('MaxAmountErrorPS')
$ws.Range('B3').Value = 'Tue Jun 03 22:34:54 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorPC')
$ws.Range('B2').Value = 'Tue Jun 03 22:38:50 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorPC')
$ws.Range('B3').Value = 'Tue Jun 03 22:39:34 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorCC')
$ws.Range('B2').Value = 'Tue Jun 03 22:35:36 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorCC')
$ws.Range('B3').Value = 'Tue Jun 03 22:36:21 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorCorp')
$ws.Range('B2').Value = 'Tue Jun 03 22:37:11 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorCorp')
$ws.Range('B3').Value = 'Tue Jun 03 22:38:05 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorPS')
$ws.Range('B2').Value = 'Tue Jun 03 22:40:17 IST 2025'

$ws = $wb.Worksheets.Item('MinAmountErrorPS')
$ws.Range('B3').Value = 'Tue Jun 03 22:41:03 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCredit')
$ws.Range('B2').Value = 'Tue Jun 03 20:47:22 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCredit')
$ws.Range('B3').Value = 'Tue Jun 03 20:48:08 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCredit')
$ws.Range('B4').Value = 'Tue Jun 03 20:48:57 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCredit')
$ws.Range('B5').Value = 'Tue Jun 03 20:49:41 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPC')
$ws.Range('B2').Value = 'Tue Jun 03 20:50:29 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPC')
$ws.Range('B3').Value = 'Tue Jun 03 20:51:12 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPC')
$ws.Range('B4').Value = 'Tue Jun 03 20:52:03 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPC')
$ws.Range('B5').Value = 'Tue Jun 03 20:52:48 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPS')
$ws.Range('B2').Value = 'Tue Jun 03 20:53:31 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPS')
$ws.Range('B3').Value = 'Tue Jun 03 20:54:21 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPS')
$ws.Range('B4').Value = 'Tue Jun 03 20:55:04 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayPS')
$ws.Range('B5').Value = 'Tue Jun 03 20:55:50 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCorp')
$ws.Range('B2').Value = 'Tue Jun 03 20:44:02 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCorp')
$ws.Range('B3').Value = 'Tue Jun 03 20:44:51 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCorp')
$ws.Range('B4').Value = 'Tue Jun 03 20:45:37 IST 2025'

$ws = $wb.Worksheets.Item('OverAndUnderPayCorp')
$ws.Range('B5').Value = 'Tue Jun 03 20:46:26 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountCorp')
$ws.Range('B2').Value = 'Tue Jun 03 22:54:03 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountCorp')
$ws.Range('B3').Value = 'Tue Jun 03 22:54:47 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountPC')
$ws.Range('B2').Value = 'Tue Jun 03 22:55:37 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountPC')
$ws.Range('B3').Value = 'Tue Jun 03 22:56:20 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountPS')
$ws.Range('B2').Value = 'Tue Jun 03 22:57:05 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountPS')
$ws.Range('B3').Value = 'Tue Jun 03 22:57:52 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountCC')
$ws.Range('B2').Value = 'Tue Jun 03 22:52:26 IST 2025'

$ws = $wb.Worksheets.Item('NoModifyAmountCC')
$ws.Range('B3').Value = 'Tue Jun 03 22:53:14 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorCC')
$ws.Range('B2').Value = 'Tue Jun 03 22:43:22 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorCC')
$ws.Range('B3').Value = 'Tue Jun 03 22:44:15 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorPC')
$ws.Range('B2').Value = 'Tue Jun 03 22:46:35 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorPC')
$ws.Range('B3').Value = 'Tue Jun 03 22:47:20 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorCorp')
$ws.Range('B2').Value = 'Tue Jun 03 22:45:03 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorCorp')
$ws.Range('B3').Value = 'Tue Jun 03 22:45:47 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorPS')
$ws.Range('B2').Value = 'Tue Jun 03 22:48:01 IST 2025'

$ws = $wb.Worksheets.Item('NoOverPayErrorPS')
$ws.Range('B3').Value = 'Tue Jun 03 22:48:43 IST 2025'

$ws = $wb.Worksheets.Item('NoUnderPayErrorPS')
$ws.Range('B2').Value = 'Tue Jun 03 22:51:44 IST 2025'

$ws = $wb.Worksheets.Item('NoUnderPayErrorPC')
$ws.Range('B2').Value = 'Tue Jun 03 22:51:00 IST 2025'

$ws = $wb.Worksheets.Item('NoUnderPayErrorCC')
$ws.Range('B2').Value = 'Tue Jun 03 22:49:27 IST 2025'

$ws = $wb.Worksheets.Item('NoUnderPayErrorCorp')
$ws.Range('B2').Value = 'Tue Jun 03 22:50:12 IST 2025'

$ws = $wb.Worksheets.Item('CardExpiredErrorCC')
$ws.Range('B2').Value = 'Tue Jun 03 22:27:24 IST 2025'

$ws = $wb.Worksheets.Item('CardExpiredErrorCC')
$ws.Range('B3').Value = 'Tue Jun 03 22:28:10 IST 2025'

$ws = $wb.Worksheets.Item('CardNotAcceptedErrorCC')
$ws.Range('B2').Value = 'Tue Jun 03 22:28:54 IST 2025'

$ws = $wb.Worksheets.Item('CardNotAcceptedErrorCC')
$ws.Range('B3').Value = 'Tue Jun 03 22:29:36 IST 2025'
